$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (and the one column-B correction at row 182) per the
# re-annotated classification values.
$ws.Range("C10").Value = 1
$ws.Range("C15").Value = 2
$ws.Range("C17").Value = 2
$ws.Range("C19").Value = 2
$ws.Range("C22").Value = 3
$ws.Range("C24").Value = 3
$ws.Range("C25").Value = 1
$ws.Range("C28").Value = 2
$ws.Range("C30").Value = 2
$ws.Range("C34").Value = 3
$ws.Range("C36").Value = 3
$ws.Range("C43").Value = 1
$ws.Range("C46").Value = 2
$ws.Range("C47").Value = 3
$ws.Range("C48").Value = 1
$ws.Range("C49").Value = 1
$ws.Range("C50").Value = 2
$ws.Range("C51").Value = 3
$ws.Range("C55").Value = 2
$ws.Range("C57").Value = 1
$ws.Range("C58").Value = 3
$ws.Range("C59").Value = 3
$ws.Range("C64").Value = 3
$ws.Range("C67").Value = 3
$ws.Range("C68").Value = 2
$ws.Range("C70").Value = 2
$ws.Range("C71").Value = 1
$ws.Range("C80").Value = 1
$ws.Range("C81").Value = 1
$ws.Range("C84").Value = 2
$ws.Range("C86").Value = 3
$ws.Range("C87").Value = 1
$ws.Range("C92").Value = 2
$ws.Range("C94").Value = 2
$ws.Range("C95").Value = 2
$ws.Range("C96").Value = 1
$ws.Range("C101").Value = 3
$ws.Range("C102").Value = 2
$ws.Range("C107").Value = 2
$ws.Range("C111").Value = 3
$ws.Range("C112").Value = 2
$ws.Range("C117").Value = 3
$ws.Range("C118").Value = 3
$ws.Range("C119").Value = 1
$ws.Range("C122").Value = 1
$ws.Range("C123").Value = 2
$ws.Range("C126").Value = 1
$ws.Range("C127").Value = 2
$ws.Range("C131").Value = 2
$ws.Range("C135").Value = 3
$ws.Range("C137").Value = 1
$ws.Range("C139").Value = 3
$ws.Range("C145").Value = 1
$ws.Range("C147").Value = 3
$ws.Range("C148").Value = 2
$ws.Range("C154").Value = 3
$ws.Range("C155").Value = 3
$ws.Range("C163").Value = 2
$ws.Range("C168").Value = 1
$ws.Range("C176").Value = 3
$ws.Range("C177").Value = 2
$ws.Range("C178").Value = 2
$ws.Range("C180").Value = 1
$ws.Range("B182").Value = 0
$ws.Range("C186").Value = 2
$ws.Range("C187").Value = 2
$ws.Range("C188").Value = 1
$ws.Range("C189").Value = 1
$ws.Range("C193").Value = 2
$ws.Range("C198").Value = 2
$ws.Range("C201").Value = 1

# Row 182 no longer carries a "Positividade" (C) annotation -- clear it so
# the cell is removed rather than left holding a stale 0.
$ws.Range("C182").ClearContents()

# Restore the view to where the editor left off scrolling/selecting.
$ws.Range("C121").Select()
